# Usman - New users added to cobaltUsers.xls
# Adds 33 new rows (81-113) of user data to the "Users" worksheet, matching
# the KHPadd / SearchWhatsMarket / SearchKnowHow / Ask / AssetPage test users.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New Name/Email pairs, in row order starting at row 81.
$newUsers = @(
    @("KHPaddUser1", "KHPaddUser1@mailinator.com"),
    @("KHPaddUser2", "KHPaddUser2@mailinator.com"),
    @("KHPaddUser3", "KHPaddUser3@mailinator.com"),
    @("KHPaddUser4", "KHPaddUser4@mailinator.com"),
    @("KHPaddUser5", "KHPaddUser5@mailinator.com"),
    @("KHPaddUser6", "KHPaddUser6@mailinator.com"),
    @("SearchWhatsMarketUser1", "SearchWhatsMarketUser1@mailinator.com "),
    @("SearchWhatsMarketUser2", "SearchWhatsMarketUser2@mailinator.com "),
    @("SearchWhatsMarketUser3", "SearchWhatsMarketUser3@mailinator.com "),
    @("SearchWhatsMarketUser4", "SearchWhatsMarketUser4@mailinator.com "),
    @("SearchWhatsMarketUser5", "SearchWhatsMarketUser5@mailinator.com "),
    @("SearchWhatsMarketUser6", "SearchWhatsMarketUser6@mailinator.com "),
    @("SearchWhatsMarketUser7", "SearchWhatsMarketUser7@mailinator.com "),
    @("SearchWhatsMarketUser8", "SearchWhatsMarketUser8@mailinator.com "),
    @("SearchKnowHowUser1", "SearchKnowHowUser1@mailinator.com "),
    @("SearchKnowHowUser2", "SearchKnowHowUser2@mailinator.com "),
    @("SearchKnowHowUser3", "SearchKnowHowUser3@mailinator.com "),
    @("SearchKnowHowUser4", "SearchKnowHowUser4@mailinator.com "),
    @("SearchKnowHowUser5", "SearchKnowHowUser5@mailinator.com "),
    @("SearchKnowHowUser6", "SearchKnowHowUser6@mailinator.com "),
    @("SearchKnowHowUser7", "SearchKnowHowUser7@mailinator.com "),
    @("SearchKnowHowUser8", "SearchKnowHowUser8@mailinator.com "),
    @("AskUser1", "AskUser1@mailinator.com "),
    @("AskUser2", "AskUser2@mailinator.com "),
    @("AskUser3", "AskUser3@mailinator.com "),
    @("AskUser4", "AskUser4@mailinator.com "),
    @("AskUser5", "AskUser5@mailinator.com "),
    @("AskUser6", "AskUser6@mailinator.com "),
    @("AssetPageUser1", "AssetPageUser1@mailinator.com "),
    @("AssetPageUser2", "AssetPageUser2@mailinator.com "),
    @("AssetPageUser3", "AssetPageUser3@mailinator.com "),
    @("AssetPageUser4", "AssetPageUser4@mailinator.com "),
    @("AssetPageUser5", "AssetPageUser5@mailinator.com ")
)

$startRow = 81
$lastRow = $startRow + $newUsers.Count - 1   # 113

# Reference cells whose existing styles we want to replicate onto the new rows.
$refRow = 80
$refB = $ws.Cells.Item($refRow, 2)
$refE = $ws.Cells.Item($refRow, 5)
$refF = $ws.Cells.Item($refRow, 6)
$refG = $ws.Cells.Item($refRow, 7)

# Build the new "Arial 10" font style by formatting A81 first; later rows copy
# this same style from A81 so only a single new style entry is created.
$ws.Cells.Item($startRow, 1).Value = $newUsers[0][0]
$ws.Cells.Item($startRow, 1).Font.Name = "Arial"
$ws.Cells.Item($startRow, 1).Font.Size = 10
$refA = $ws.Cells.Item($startRow, 1)

for ($i = 0; $i -lt $newUsers.Count; $i++) {
    $row = $startRow + $i
    $name = $newUsers[$i][0]
    $email = $newUsers[$i][1]

    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)
    $cellE = $ws.Cells.Item($row, 5)
    $cellF = $ws.Cells.Item($row, 6)
    $cellG = $ws.Cells.Item($row, 7)

    if ($row -ne $startRow) {
        $refA.Copy($cellA)
    }
    $cellA.Value = $name

    $refB.Copy($cellB)
    $cellB.Value = "Password1"

    $refE.Copy($cellE)
    $cellE.Value = "THIS IS IN USE 24/7 - DO NOT USE!"

    $refF.Copy($cellF)
    $cellF.Value = "N"

    $refG.Copy($cellG)
    $cellG.Value = $email
}

# Hyperlinks on column G - added in this specific order (matches the source
# workbook); rows 81 and 112 are intentionally left without hyperlinks, and
# row 111's hyperlink carries a stale display tooltip from AssetPageUser1.
# Adding a hyperlink forces Excel to re-apply cell formatting, so the
# original "s=10" style (copied from G80) is restored afterwards to keep the
# new cells visually consistent with the rest of column G.
$hyperlinkRows = @(82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,113)

foreach ($row in $hyperlinkRows) {
    $email = ($newUsers[$row - $startRow][1]).Trim()
    $cellG = $ws.Cells.Item($row, 7)
    $ws.Hyperlinks.Add($cellG, "mailto:$email", [Type]::Missing, [Type]::Missing, [Type]::Missing)
    $refG.Copy($cellG)
    $cellG.Value = $newUsers[$row - $startRow][1]
}

# Row 111 (AssetPageUser3) keeps its normal text but the hyperlink's display
# tooltip mistakenly shows AssetPageUser1's address, as in the source data.
$cellG111 = $ws.Cells.Item(111, 7)
$ws.Hyperlinks.Add($cellG111, "mailto:AssetPageUser3@mailinator.com", [Type]::Missing, [Type]::Missing, "AssetPageUser1@mailinator.com ")
$refG.Copy($cellG111)
$cellG111.Value = "AssetPageUser3@mailinator.com "

# Update the visible selection to span the newly added rows.
$selected = $ws.Range("A81:G113").Select()

Write-Host "Added rows $startRow to $lastRow with hyperlinks"
